# Update the "Förändrad" (Changed) date column (C) from serial 45599 (2024-11-03)
# to serial 45600 (2024-11-04) for all data rows (rows 2 through 32).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 32; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45599) {
        $cell.Value2 = 45600
    }
}
